# Add a web parser class and function to work with the GPT API --
# in terms of document content this appends four new paragraphs
# (an in-text citation sentence followed by its reference entry,
# repeated for citations [1] and [2]) at the very end of the document,
# mirroring the paragraphs already present just above them.

$d = $word.ActiveDocument

$newParagraphs = @(
    'Scientific writing is a cornerstone of scholarly communication, facilitating the dissemination of knowledge, discoveries, and advancements across various disciplines. [1].',
    '[1]  Berco A, "DSL: Scientific text processing". TUM. 2010 Available from: www.overleaf/PBLTeam_6. Accesed Date 02 May 2024.',
    'Scientific writing encompasses various genres, including research papers, reviews, and technical reports [2].',
    '[2]  Berco A, "DSL: Scientific text processing". TUM. 2010 Available from: www.overleaf/PBLTeam_6. Accesed Date 02 May 2024.'
)

foreach ($t in $newParagraphs) {
    $tail = $d.Content
    $tail.Collapse(0)
    $tail.InsertParagraphAfter()
    $tail.Collapse(0)
    $tail.Text = $t
}
